$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:207 down to 104:208.
$ws.Rows("103").Insert()

# Populate the new row 103 with the new data point.
$ws.Range("A103").Value = 8
$ws.Range("B103").Value = "Terminal La Palmera de La Serena"
$ws.Range("C103").Value = "Coquimbo"
$ws.Range("D103").Value = 45264
$ws.Range("E103").Value = 4
$ws.Range("F103").Value = 100114007
$ws.Range("G103").Value = "Jengibre"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 440
$ws.Range("K103").Value = 22000
$ws.Range("L103").Value = 23000
$ws.Range("M103").Value = 22500
$ws.Range("N103").Value = "$/caja 13 kilos"
$ws.Range("O103").Value = "Perú"
$ws.Range("P103").Value = 1731
$ws.Range("Q103").Value = 13
$ws.Range("R103").Value = "Hortaliza"
